$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.013.80'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '1.910.10'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''0.7906'
$ws.Range("E5").Value = '  +4.58%  '
$ws.Range("D6").Value = '''242.04'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '''0.3161'
$ws.Range("E8").Value = '  +2.69%  '
$ws.Range("D9").Value = '''26.31'
$ws.Range("E9").Value = '  +3.00%  '
$ws.Range("D10").Value = '''0.06898'
$ws.Range("E10").Value = '  -0.13%  '
$ws.Range("D11").Value = '''0.08002'
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.909.09'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '''0.7428'
$ws.Range("E13").Value = '  -1.74%  '
$ws.Range("D14").Value = '''5.186'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").Value = '''93.04'
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Value = '30.009.61'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").Value = '''5.864'
$ws.Range("E18").Value = '  -5.34%  '
$ws.Range("D19").Value = '''245.61'
$ws.Range("E19").Value = '  +3.27%  '
$ws.Range("D20").Value = '''0.000007737'
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = '2.153.09'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").Value = '''1.001'
$ws.Range("D24").Value = '''6.832'
$ws.Range("E24").Value = '  -3.12%  '
$ws.Range("D25").Value = '''168.01'
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("D26").Value = '''9.224'
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("D27").Value = '''0.1387'
$ws.Range("E27").Value = '  +7.61%  '
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("D29").Value = '''2.031'
$ws.Range("E29").Value = '  -2.01%  '
$ws.Range("D30").Value = '''1.365'
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("D32").Value = '''4.313'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("E33").Value = '  +0.81%  '
$ws.Range("D34").Value = '''0.05519'
$ws.Range("E34").Value = '  +2.30%  '
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("D36").Value = '''0.7324'
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("D37").Value = '''2.722'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '''0.01924'
$ws.Range("E38").Value = '  -1.17%  '
$ws.Range("D39").Value = '''2.784'
$ws.Range("E39").Value = '  +0.69%  '
$ws.Range("D40").Value = '''6.131'
$ws.Range("E40").Value = '  -1.97%  '
$ws.Range("D41").Value = '''0.4410'
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("D42").Value = '''72.37'
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '''0.8374'
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("D45").Value = '''1.873'
$ws.Range("E45").Value = '  -3.91%  '
$ws.Range("D46").Value = '''100.38'
$ws.Range("E46").Value = '  -1.23%  '
$ws.Range("D47").Value = '''7.546'
$ws.Range("E47").Value = '  -2.13%  '
$ws.Range("D48").Value = '''987.78'
$ws.Range("E48").Value = '  +7.38%  '
$ws.Range("D49").Value = '2.056.49'
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").Value = '''36.28'
$ws.Range("E50").Value = '  -0.79%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.05937'
$ws.Range("E51").Value = '  -0.73%  '
